# data-siswa.xlsx : replace the PKL student roster with the final list.
#
# The NIM column holds values that look numeric ("151080200212", ...).
# A plain  $range.Value = "151080200212"  would make Excel infer a
# genuine Number for that cell (losing the text/shared-string typing the
# workbook originally used). To keep every cell a proper text value -
# exactly like the original file - we stage each value in a scratch cell
# with a leading apostrophe (forcing Excel to treat it as text), copy
# that cell, and Paste-Special "Values only" into the real destination.
# Paste Special (values only) carries over the source's text typing
# without touching the destination cell's own formatting, so no stray
# NumberFormat/style ends up attached to the edited cells - just like
# the diff we're replicating, which only touches the shared string
# table content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z99")

function Set-TextValue {
    param($range, [string]$text)

    $scratch.Value = "'" + $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Delete()
}

# NIM / SISWA rows 2-5, in final order.
Set-TextValue $ws.Range("A2") "1510800200225"
Set-TextValue $ws.Range("B2") "APRILIA KUSUMA NINGRUM"

Set-TextValue $ws.Range("A3") "151080200212"
Set-TextValue $ws.Range("B3") "AHBABUL MUSTHOFA"

Set-TextValue $ws.Range("A4") "151080200224"
Set-TextValue $ws.Range("B4") "M.AGUNG HERIYANTO"

Set-TextValue $ws.Range("A5") "191080200223"
Set-TextValue $ws.Range("B5") "M.AUNUR ROSIDIN"
